# chart_config.xlsx — "multiple definitions working now"
#
# Adds a second panel ("Gasoil-Heating oil box spread") alongside the
# existing Brent-WTI one, wired through all three tables:
#   panels      (sheet1) - new "extend_tenor"/"extend_count" columns + new panel row
#   products    (sheet2) - new product row (mt-denominated) for the new panel
#   expressions (sheet3) - two new expression rows (GO-F, HO-F) for the new panel

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# products — add product_id 2 (panel_id 2, USD, mt, x, y1)
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("products")
$wsProducts.Range("A3").Value = 2
$wsProducts.Range("B3").Value = 2
$wsProducts.Range("C3").Value = "USD"
$wsProducts.Range("D3").Value = "mt"
$wsProducts.Range("E3").Value = "x"
$wsProducts.Range("F3").Value = "y1"

# ---------------------------------------------------------------------
# expressions — add expression_id 3 (GO-F) and 4 (HO-F) for product_id 2
# ---------------------------------------------------------------------
$wsExpr = $wb.Worksheets.Item("expressions")
$wsExpr.Range("A4").Value = 3
$wsExpr.Range("B4").Value = 2
$wsExpr.Range("C4").Value = 1
$wsExpr.Range("D4").Value = "GO-F"
$wsExpr.Range("E4").Value = 202203
$wsExpr.Range("G4").Value = 202204

$wsExpr.Range("A5").Value = 4
$wsExpr.Range("B5").Value = 2
$wsExpr.Range("C5").Value = -1
$wsExpr.Range("D5").Value = "HO-F"
$wsExpr.Range("E5").Value = 202203
$wsExpr.Range("G5").Value = 202204

# ---------------------------------------------------------------------
# panels — add "extend_tenor" / "extend_count" columns for the existing
# panel, then the new panel_id 2 row (Gasoil-Heating oil box spread)
# ---------------------------------------------------------------------
$wsPanels = $wb.Worksheets.Item("panels")
$wsPanels.Range("D1").Value = "extend_tenor"
$wsPanels.Range("D2").Value = "month"
$wsPanels.Range("D3").Value = "month"
$wsPanels.Range("E1").Value = "extend_count"
$wsPanels.Range("E2").Value = 3
$wsPanels.Range("E3").Value = 6
$wsPanels.Range("B3").Value = "Gasoil-Heating oil  box spread {front}-{back}"
$wsPanels.Range("A3").Value = 2
$wsPanels.Range("C3").Value = 3

# ---------------------------------------------------------------------
# Restore the per-sheet selections, finishing on "expressions" so it
# stays the active tab.
# ---------------------------------------------------------------------
$wsPanels.Range("B4").Select() | Out-Null
$wsProducts.Range("G21").Select() | Out-Null
$wsExpr.Range("I21").Select() | Out-Null
